$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that follows the title.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold "Play Blazin Hot 7s Stack Em Up for Free - Game Review"
#    paragraph right before the final (italic) paragraph.
$n = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs.Item($n - 1)
$insertPos = $secondToLast.Range.End - 1
$insertionPoint = $d.Range($insertPos, $insertPos)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Blazin Hot 7s Stack Em Up for Free - Game Review</w:t></w:r></w:p>'
$insertionPoint.InsertXML($newParaXml)

# 3. Replace the old "Create a feature image..." prompt text in the final
#    paragraph with the meta-description copy (keeping its italic run).
$d.Content.Find.Execute(
    'Create a feature image fitting "Blazin Hot 7s Stack Em Up": - Draw a cartoon-style image of a happy Maya warrior with glasses wearing a headdress made of fruits such as cherries, oranges, lemons, plums, and watermelons. - Have the warrior holding a Stack''Em Up symbol in one hand and a handful of coins in the other hand. - Surround the warrior with cascading reels and colorful symbols. - Add text above the image that says "Blazin Hot 7s Stack Em Up" in bold, fiery letters.',
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discover the unique game mechanics and bonus features of Blazin Hot 7s Stack Em Up. Play for free and win big with high volatility and RTP.",
    2)
